$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Settings")
$ws2 = $wb.Worksheets.Item("Constants")

$ws1.Range("B46").Select()
$ws2.Activate()
$ws2.Range("A17").Select()
